# Funcional (aunque desmodalurizada) version del programa
#
# - Clears the old "dgdgh"/"sss" values out of column E (rows 1-6).
# - Leaves E2 with an (empty) underlined-font style, matching the
#   underline formatting that was left behind on that cell.
# - Fills in six new names under column D (rows 6-11): Ramon, Ramon,
#   Pedro, Pedro, Pedro, Marcos.
# - Gives the last new entry (D11) its own distinct cell style.
# - Selects the whole D column, matching the final selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old column-E values (the "dgdgh" / "sss" shared strings
# disappear from the workbook once nothing references them any more).
$ws.Range("E1:E6").ClearContents()

# E2 keeps an underlined-font style even though it is now blank.
$ws.Range("E2").Font.Underline = $true

# New names appended under column D.
$ws.Range("D6").Value = "Ramon"
$ws.Range("D7").Value = "Ramon"
$ws.Range("D8").Value = "Pedro"
$ws.Range("D9").Value = "Pedro"
$ws.Range("D10").Value = "Pedro"
$ws.Range("D11").Value = "Marcos"

# D11 carries its own (visually default) cell style, distinct from the
# plain unstyled cells above it.
$ws.Range("D11").Font.ThemeColor = 1

# Final selection: the whole of column D.
$ws.Columns("D:D").Select() | Out-Null
